$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells, matching formatting of existing header cells
$ws.Range("G1").Value = "Elapsed Time"
$ws.Range("H1").Value = "CPU"

# Copy formatting (font, fill, border, alignment) from F1 into G1:H1
$ws.Range("F1").Copy() | Out-Null
$ws.Range("G1:H1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# Update existing row 2 values
$ws.Range("B2").Value = 0.1536395509049464
$ws.Range("C2").Value = 0.9984164699088494
$ws.Range("D2").Value = 0.3311439993263866
$ws.Range("F2").Value = "Pipeline(steps=[('model', AdaBoostRegressor(learning_rate=0.5))])"

# Add new values for the added columns
$ws.Range("G2").Value = 0.1194315095165318
$ws.Range("H2").Value = 0.9890000000000001
